$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

try {
    $s1.ApplyThemeColorScheme("dummy")
    Write-Output "ApplyThemeColorScheme OK"
} catch {
    Write-Output ("ApplyThemeColorScheme ERR: " + $_.Exception.Message)
}

try {
    $s1.ApplyTemplate("dummy.thmx")
    Write-Output "ApplyTemplate OK"
} catch {
    Write-Output ("ApplyTemplate ERR: " + $_.Exception.Message)
}

try {
    $p.ApplyTemplate("dummy.thmx")
    Write-Output "Presentation.ApplyTemplate OK"
} catch {
    Write-Output ("Presentation.ApplyTemplate ERR: " + $_.Exception.Message)
}
